# Apply "All Country Files Saved And Formatted" edit to the Dominica sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fix header text + add new columns of data (G = TotalConfirmedNewCases,
#    I = TotalNewDeaths) for the existing rows, and append a new data row (9).
# ---------------------------------------------------------------------------

$ws.Range("L1").Value = "MasterSheet RowNo."

$ws.Range("G2").Value = 1
$ws.Range("I2").Value = 0

$ws.Range("G3").Value = 1
$ws.Range("I3").Value = 0

$ws.Range("G4").Value = 5
$ws.Range("I4").Value = 0

$ws.Range("G5").Value = 4
$ws.Range("I5").Value = 0

$ws.Range("G6").Value = 0
$ws.Range("I6").Value = 0

$ws.Range("G7").Value = 0
$ws.Range("I7").Value = 0

$ws.Range("G8").Value = 0
$ws.Range("I8").Value = 0

# New row 9
$ws.Range("A9").Value = 71
$ws.Range("B9").Value = 57
$ws.Range("C9").Value = "LATIN AMER. & CARIB    "
$ws.Range("D9").Value = 43921
$ws.Range("E9").Value = "Dominica"
$ws.Range("F9").Value = 11
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = "Local transmission"
$ws.Range("K9").Value = 4
$ws.Range("L9").Value = 5379

# ---------------------------------------------------------------------------
# 2. Column widths: A:O = 27
# ---------------------------------------------------------------------------
for ($c = 1; $c -le 15; $c++) {
  $ws.Columns.Item($c).ColumnWidth = 26.1666666666667
}

# ---------------------------------------------------------------------------
# 3. Styling: center horizontal/vertical alignment across A1:O9, and a
#    distinct date format ("yyyy-mm-dd;") + centered alignment for column D.
#    A seed cell is fully formatted first, then its format is propagated via
#    Copy + PasteSpecial(formats) so the runtime doesn't fabricate extra
#    intermediate cell-style entries.
# ---------------------------------------------------------------------------

# General style seed: A1
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").VerticalAlignment = -4108

# Date style seed: D1
$ws.Range("D1").NumberFormat = "yyyy-mm-dd;"
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("D1").VerticalAlignment = -4108

# Propagate the general style to all cells except column D
$ws.Range("A1").Copy()
$ws.Range("A1:C9").PasteSpecial(-4122)
$ws.Range("E1:O9").PasteSpecial(-4122)

# Propagate the date style to the rest of column D
$ws.Range("D1").Copy()
$ws.Range("D2:D9").PasteSpecial(-4122)

$excel.CutCopyMode = 0
